$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.19090754144846 / 100000
$ws.Range("C2").Value = 0.002658071450198252
$ws.Range("D2").Value = 18.71679738969934
$ws.Range("E2").Value = 2459690191846.092
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 2459690191864.811
